$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift all timestamps in column A (rows 2-97) forward by 28 days
for ($r = 2; $r -le 97; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $cell.Value = $cell.Value2 + 28
}

# Update Actual Production (MW) values in column B for rows 2-42
$bValues = @{
    2 = 1730
    3 = 1777
    4 = 1834
    5 = 1872
    6 = 1947
    7 = 1991
    8 = 2040
    9 = 2049
    10 = 2094
    11 = 2098
    12 = 2099
    13 = 2104
    14 = 2096
    15 = 2095
    16 = 2069
    17 = 2083
    18 = 2089
    19 = 2084
    20 = 2028
    21 = 2060
    22 = 2169
    23 = 2155
    24 = 2135
    25 = 2157
    26 = 2166
    27 = 2164
    28 = 2189
    29 = 2126
    30 = 2151
    31 = 2196
    32 = 2245
    33 = 2255
    34 = 2231
    35 = 2292
    36 = 2327
    37 = 2346
    38 = 2391
    39 = 2334
    40 = 2350
    41 = 2372
    42 = 2377
}
foreach ($r in $bValues.Keys) {
    $ws.Cells.Item([int]$r, 2).Value = $bValues[$r]
}
